$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 8: E8 changes from "Pendiente ADM" to "02013965 " ---
# Force text formatting so the leading zero and trailing space survive.
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "02013965 "

# --- Add new row 9 ---
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "7146"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "12/29/2025"
$ws.Range("C9").Value = "TUCUMAN 1511"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Pendiente ADM"
$ws.Range("F9").Value = "Optical Power"
$ws.Range("G9").Value = "Pendiente"
$ws.Range("H9").Value = "tendido bajo"
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = '{"direccionesNormalizadas": [{"altura": 1511, "cod_calle": 21060, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.388356", "y": "-34.601692"}, "direccion": "TUCUMAN 1511, CABA", "nombre_calle": "TUCUMAN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K9").Value = -58.388356
$ws.Range("L9").Value = -34.601692
$ws.Range("M9").Value = "San Telmo"
$ws.Range("N9").Value = "Capital Sur"

# --- Add new row 10 ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "7987"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "12/30/2025"
$ws.Range("C10").Value = "Terrada 2309"
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = "Pendiente ADM"
$ws.Range("F10").Value = "Optical Power"
$ws.Range("G10").Value = "Pendiente"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "cables colgando "
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = '{"direccionesNormalizadas": [{"altura": 2309, "cod_calle": 21021, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.482084", "y": "-34.608289"}, "direccion": "TERRADA 2309, CABA", "nombre_calle": "TERRADA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K10").Value = -58.482084
$ws.Range("L10").Value = -34.608289
$ws.Range("M10").Value = "Paternal"
$ws.Range("N10").Value = "Capital Norte"
